# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml, the
# theme wired to the one-and-only slide master) while ppt/theme/theme1.xml
# (wired to the notes master) holds the stock "Office Theme" palette.
#
# The target edit swaps the two: the slide master's theme becomes the
# standard Office Theme color palette, Integral's colors move to the
# notes-master-side theme part. The font scheme and format scheme are
# identical between the two themes (only the 12 color-scheme entries and
# the theme/clrScheme display names differ), so the visible effect of the
# swap is fully captured by re-pointing the 12 scheme colors that PowerPoint
# exposes through Slide.ThemeColorScheme (shared by every slide, since the
# deck has a single slide master/design).

function Hex2RGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order exposed by ThemeColorScheme.Item(n):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = Hex2RGB($officeTheme[$i - 1])
}

Write-Host "Applied Office Theme color scheme to the deck's theme (was Integral)."
